$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sierra Leone master data: language code column changes from French ("fra")
# to English ("eng") for every authentication-method row.
$ws.Range("A2:A6").Value = "eng"
